$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data: MaSoThue (column G) ---
$ws.Range("G2").Value = 1234567890
$ws.Range("G3").Value = 1234567891
$ws.Range("G4").Value = 1234567892
$ws.Range("G5").Value = 1234567893
$ws.Range("G6").Value = 1234567894
$ws.Range("G7").Value = 1234567895

# --- New data: MaSoHopDong (column J) ---
$ws.Range("J2").Value = "HD2468013579"
$ws.Range("J3").Value = "HD2468013580"
$ws.Range("J4").Value = "HD2468013581"
$ws.Range("J5").Value = "HD2468013582"
$ws.Range("J6").Value = "HD2468013583"
$ws.Range("J7").Value = "HD2468013584"

# --- New data: MaCongTo (column L) ---
$ws.Range("L2").Value = "ABCXYZ"
$ws.Range("L3").Value = "ABCXYZ"
$ws.Range("L4").Value = "ABCXYZ"
$ws.Range("L5").Value = "ABCXYZ"
$ws.Range("L6").Value = "ABCXYZ"
$ws.Range("L7").Value = "ABCXYZ"

# --- New data: SoNganHang (column M) - stored as text via leading quote ---
$ws.Range("M2").Value = "'5907000000001"
$ws.Range("M3").Value = "'5907000000002"
$ws.Range("M4").Value = "'5907000000003"
$ws.Range("M5").Value = "'5907000000004"
$ws.Range("M6").Value = "'5907000000005"
$ws.Range("M7").Value = "'5907000000006"

# --- New data: TenNganHang (column N) ---
$ws.Range("N2").Value = "Vietcombank"
$ws.Range("N3").Value = "Vietcombank"
$ws.Range("N4").Value = "Vietcombank"
$ws.Range("N5").Value = "Vietcombank"
$ws.Range("N6").Value = "Vietcombank"
$ws.Range("N7").Value = "Vietcombank"

# --- Column width adjustments (fix column widths for DataGridView display) ---
$ws.Columns("G").ColumnWidth = 10.166666666666666
$ws.Columns("M").ColumnWidth = 14.8776041666666

# --- Selection change ---
$ws.Range("O6").Select()
